# Update the build/version timestamp throughout the workbook.
# Old:  mines - January 30 (built on January 30 2026 16.19.47 EST)
# New:  mines - January 30 (built on February 02 2026 12.49.33 EST)

$newVersion = "mines - January 30 (built on February 02 2026 12.49.33 EST)"

$wb = $excel.ActiveWorkbook

# --- "About" sheet ---
$about = $wb.Worksheets.Item("About")

$about.Range("A2").Value = "Version: $newVersion"

$about.Range("A6").Value = "Recommended Citation:  `"Global Energy Monitor, Coal mine boundaries and methane sources for Oaky Creek Coal Mine, Australia, M0085, version '$newVersion'. (See the CC license for attribution requirements if sharing or adapting the data set.)"

# --- "Boundaries and methane sources" sheet ---
$data = $wb.Worksheets.Item("Boundaries and methane sources")

# Column S holds the build_version value for every data row (rows 2-66).
$data.Range("S2:S66").Value = $newVersion
